$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Ajo" (garlic) sheet.
# It belongs right after the current row 282 (chronologically it slots in
# as the new row 283), so insert a fresh row there which pushes the
# existing rows 283:334 down to 284:335.
$ws.Rows("283:283").Insert()

$ws.Range("A283").Value = 11
$ws.Range("B283").Value = 'Vega Monumental Concepción'
$ws.Range("C283").Value = 'Bíobío'
$ws.Range("D283").Value = 45211
$ws.Range("E283").Value = 8
$ws.Range("F283").Value = 100112003
$ws.Range("G283").Value = 'Ajo'
$ws.Range("H283").Value = 'Chino'
$ws.Range("I283").Value = 'Primera'
$ws.Range("J283").Value = 200
$ws.Range("K283").Value = 19000
$ws.Range("L283").Value = 20000
$ws.Range("M283").Value = 19500
$ws.Range("N283").Value = '$/caja 10 kilos'
$ws.Range("O283").Value = 'China'
$ws.Range("P283").Value = 1950
$ws.Range("Q283").Value = 10
$ws.Range("R283").Value = 'Hortaliza'
